# Spring 24 week 13 inputs — append new matchup rows to the "Nine" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(1, 18, 3, 2),
    @(5, 20, 6, 0),
    @(5, 13, 4, 7),
    @(4, 6, 3, 14),
    @(5, 12, 6, 8),
    @(4, 6, 3, 14),
    @(4, 4, 5, 16),
    @(5, 7, 8, 13),
    @(2, 12, 1, 8),
    @(5, 8, 4, 12),
    @(4, 6, 3, 14),
    @(5, 3, 7, 17),
    @(3, 0, 1, 20),
    @(4, 5, 8, 15),
    @(5, 12, 4, 8),
    @(7, 8, 5, 12),
    @(5, 15, 3, 5),
    @(2, 12, 5, 8),
    @(4, 6, 5, 14),
    @(4, 3, 3, 17),
    @(6, 12, 7, 8),
    @(7, 6, 6, 14),
    @(2, 7, 4, 13)
)

$rowCount = $data.Count
$colCount = 4

# Existing data occupies rows 1..1574 (A1:D1574); new rows start right after.
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$lastRow = $startRow + $rowCount - 1

# Build a proper 2D array so the whole block can be written in one shot.
$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$targetRange = $ws.Range("A$startRow" + ":D$lastRow")
$targetRange.Value = $values

$nextCell = "A" + ($lastRow + 1)
$ws.Range($nextCell).Select() | Out-Null
